$wb = $excel.ActiveWorkbook
Write-Output $wb.Worksheets.Count
for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Output $wb.Worksheets.Item($i).Name
}
